$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("C8").Value = "Volume 30   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/10/2023  Through  4/16/2023"

# --- Crime statistics table updates (rows 14-29) ---
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 3
$ws.Range("H14").Value = 200
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 33.333333333333
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -10
$ws.Range("L15").Value = 28.571428571428
$ws.Range("N15").Value = -35.714285714285
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 166.666666666667
$ws.Range("F16").Value = 39
$ws.Range("G16").Value = 29
$ws.Range("H16").Value = 34.482758620689
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = -18.965517241379
$ws.Range("L16").Value = 6.818181818181
$ws.Range("M16").Value = 2.173913043478
$ws.Range("N16").Value = -72.674418604651
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 5.555555555555
$ws.Range("F17").Value = 71
$ws.Range("G17").Value = 47
$ws.Range("H17").Value = 51.063829787234
$ws.Range("I17").Value = 179
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = 14.012738853503
$ws.Range("L17").Value = 54.310344827586
$ws.Range("M17").Value = 108.139534883721
$ws.Range("N17").Value = 6.547619047619
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -90.909090909090
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -19.230769230769
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = -36.363636363636
$ws.Range("L18").Value = -5.084745762711
$ws.Range("M18").Value = -5.084745762711
$ws.Range("N18").Value = -85.786802030456
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 8.823529411764
$ws.Range("I19").Value = 127
$ws.Range("J19").Value = 132
$ws.Range("K19").Value = -3.787878787878
$ws.Range("L19").Value = 22.115384615384
$ws.Range("M19").Value = 67.105263157894
$ws.Range("N19").Value = -7.299270072992
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 36.842105263157
$ws.Range("I20").Value = 101
$ws.Range("J20").Value = 78
$ws.Range("K20").Value = 29.487179487179
$ws.Range("L20").Value = 180.555555555556
$ws.Range("M20").Value = 304
$ws.Range("N20").Value = -42.937853107344
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 2.272727272727
$ws.Range("F21").Value = 199
$ws.Range("G21").Value = 158
$ws.Range("H21").Value = 25.949367088607
$ws.Range("I21").Value = 574
$ws.Range("J21").Value = 587
$ws.Range("K21").Value = -2.214650766609
$ws.Range("L21").Value = 38.313253012048
$ws.Range("M21").Value = 66.376811594202
$ws.Range("N21").Value = -53.969526864474
$ws.Range("C22").NumberFormat = "General"
$ws.Range("C22").Value = "'0"
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = -50
$ws.Range("M22").Value = -66.666666666666
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 11
$ws.Range("K23").Value = -36.363636363636
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -23.529411764705
$ws.Range("F24").Value = 66
$ws.Range("H24").Value = -33.333333333333
$ws.Range("I24").Value = 288
$ws.Range("J24").Value = 306
$ws.Range("K24").Value = -5.882352941176
$ws.Range("L24").Value = 29.729729729729
$ws.Range("M24").Value = 15.662650602409
$ws.Range("C25").Value = 20
$ws.Range("E25").Value = 11.111111111111
$ws.Range("F25").Value = 81
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = 30.645161290322
$ws.Range("I25").Value = 224
$ws.Range("J25").Value = 262
$ws.Range("K25").Value = -14.503816793893
$ws.Range("L25").Value = 17.894736842105
$ws.Range("M25").Value = -3.862660944206
$ws.Range("J26").Value = 17
$ws.Range("K26").Value = -11.764705882352
$ws.Range("L26").Value = 7.142857142857
$ws.Range("C27").NumberFormat = "General"
$ws.Range("C27").Value = "'0"
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -60
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = -21.739130434782
$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = -12.5
$ws.Range("L28").Value = 16.666666666666
$ws.Range("N28").Value = -53.333333333333
$ws.Range("C29").NumberFormat = "General"
$ws.Range("C29").Value = "'0"
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = -33.333333333333
$ws.Range("J29").Value = 14
$ws.Range("K29").Value = -42.857142857142
$ws.Range("L29").Value = -27.272727272727
$ws.Range("N29").Value = -69.230769230769
